$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (North America)
$ws.Range("B2").Value = 126947780
$ws.Range("D2").Value = 1636848
$ws.Range("F2").Value = 122930671
$ws.Range("G2").Value = 517
$ws.Range("H2").Value = 2380261
$ws.Range("I2").Value = 5953

# Row 3 (Asia)
$ws.Range("B3").Value = 218123983
$ws.Range("C3").Value = 669
$ws.Range("D3").Value = 1547511
$ws.Range("E3").Value = 4
$ws.Range("F3").Value = 201785146
$ws.Range("G3").Value = 1208
$ws.Range("H3").Value = 14791326
$ws.Range("I3").Value = 15123

# Row 4 (Europe)
$ws.Range("B4").Value = 249633231
$ws.Range("C4").Value = 257
$ws.Range("D4").Value = 2065221
$ws.Range("E4").ClearContents()
$ws.Range("F4").Value = 245715733
$ws.Range("G4").Value = 1329
$ws.Range("H4").Value = 1852277
$ws.Range("I4").Value = 5532

# Row 5 (South America)
$ws.Range("B5").Value = 68801962
$ws.Range("D5").Value = 1357193
$ws.Range("F5").Value = 66484069
$ws.Range("G5").ClearContents()
$ws.Range("H5").Value = 960700
$ws.Range("I5").Value = 10100

# Row 6 (Australia/Oceania)
$ws.Range("B6").Value = 14521437
$ws.Range("D6").Value = 28978
$ws.Range("F6").Value = 14355427
$ws.Range("H6").Value = 137032
$ws.Range("I6").Value = 86

# Row 7 (Africa)
$ws.Range("B7").Value = 12829479
$ws.Range("D7").Value = 258804
$ws.Range("F7").Value = 12087469
$ws.Range("H7").Value = 483206
$ws.Range("I7").Value = 547
